$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 25 with level 1 block entry
$ws.Range("A25").Value = 18
$ws.Range("B25").Value = "#"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "Bloc"

# Update selection to match the final state
$ws.Range("D25").Select()
